$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New log rows captured on 2018.08.23 and 2018.08.24.
# Columns: Date, Time, Neuron Type, Radius (nm), Thickness (um), Fdrive (kHz),
#          Adrive (kPa), Tstim (ms), PRF (kHz), Duty factor, Sim. Type,
#          # samples, Comp. time (s), # spikes, Latency (ms), Spike rate (sp/ms)
$data = @(
    @("2018.08.23", "18:45:03", "RS", 32, 0, 2010, 305, 250, 5.005, 0.505, "effective", 6004, 16.06, 181, 40.95904095904093, 0.8630950329170262),
    @("2018.08.23", "18:50:13", "RS", 32, 0, 500, 100, 250, 0.1, 0.9500000000000001, "effective", 6000, 9.17, 78, 37.1, 0.4857362089357296),
    @("2018.08.23", "19:04:08", "RS", 32, 0, 500, 100, 250, 0.1, 0.9, "effective", 6000, 8.039999999999999, 60, 39.09999999999999, 0.4846200285278927),
    @("2018.08.24", "16:37:39", "RS", 32, 0, 500, 100, 250, 0.1, 0.96, "effective", 6000, 9.720000000000001, 85, 36.8, 0.4878273340644669),
    @("2018.08.24", "16:39:50", "RS", 32, 0, 500, 95, 250, "N/A", 1, "effective", 6000, 10.79, 111, 37.05, 0.515805960224653),
    @("2018.08.24", "16:42:28", "RS", 32, 0, 500, 11, 250, "N/A", 1, "effective", 6000, 1.65, 0, "N/A", "N/A")
)

$startRow = 10
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]
    for ($c = 0; $c -lt $rowData.Length; $c++) {
        $col = $c + 1
        $cell = $ws.Cells.Item($row, $col)
        $value = $rowData[$c]
        # The "Date" column holds strings like "2018.08.23" that Excel would
        # otherwise auto-detect and silently convert to a date serial number.
        # Force the cell to Text format first so the literal string is kept,
        # then clear the (sticky) cell format so no extra style is left behind.
        if ($col -eq 1) {
            $cell.NumberFormat = "@"
            $cell.Value = $value
            $cell.ClearFormats()
        } else {
            $cell.Value = $value
        }
    }
}
